$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run RU 1001 (without crop) produced #NUM! errors for the
# "Saudi Arabia" column (C) on each data row.
$ws.Range("C2").Value = "#NUM!"
$ws.Range("C3").Value = "#NUM!"
$ws.Range("C4").Value = "#NUM!"
$ws.Range("C5").Value = "#NUM!"
$ws.Range("C6").Value = "#NUM!"
